{"js": "// Lattice multiplication exercises: update the numbers in each table cell.\n// Each cell holds a single paragraph/run whose text is split into 5 segments\n// joined by manual line breaks (\"\\u000b\" <-> <w:br/>):\n//   \"AB x CD\", \"  C    D\", \"  ----\", \"X|    |\", \"Y|    |\"\n// We replace, cell by cell (row-major order), the 5 segments with their new\n// values while keeping the run's own formatting (font size) intact by\n// rewriting through the paragraph's own Range rather than the cell body\n// (body-level replace would drop the run's rPr).\n\nconst newCellLines = [\n  [\"74 x 25\", \"  2    5\", \"  ----\", \"7|    |\", \"4|    |\"],\n  [\"52 x 41\", \"  4    1\", \"  ----\", \"5|    |\", \"2|    |\"],\n  [\"43 x 24\", \"  2    4\", \"  ----\", \"4|    |\", \"3|    |\"],\n  [\"90 x 42\", \"  4    2\", \"  ----\", \"9|    |\", \"0|    |\"],\n  [\"62 x 15\", \"  1    5\", \"  ----\", \"6|    |\", \"2|    |\"],\n  [\"91 x 28\", \"  2    8\", \"  ----\", \"9|    |\", \"1|    |\"],\n  [\"64 x 15\", \"  1    5\", \"  ----\", \"6|    |\", \"4|    |\"],\n  [\"61 x 61\", \"  6    1\", \"  ----\", \"6|    |\", \"1|    |\"],\n  [\"77 x 77\", \"  7    7\", \"  ----\", \"7|    |\", \"7|    |\"],\n  [\"57 x 52\", \"  5    2\", \"  ----\", \"5|    |\", \"7|    |\"],\n  [\"97 x 29\", \"  2    9\", \"  ----\", \"9|    |\", \"7|    |\"],\n  [\"37 x 65\", \"  6    5\", \"  ----\", \"3|    |\", \"7|    |\"],\n  [\"50 x 74\", \"  7    4\", \"  ----\", \"5|    |\", \"0|    |\"],\n  [\"34 x 35\", \"  3    5\", \"  ----\", \"3|    |\", \"4|    |\"],\n  [\"17 x 16\", \"  1    6\", \"  ----\", \"1|    |\", \"7|    |\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst rows = table.rowCount;\nconst cols = 3;\n\nlet cellIndex = 0;\nfor (let r = 0; r < rows; r++) {\n  for (let c = 0; c < cols; c++) {\n    if (cellIndex >= newCellLines.length) break;\n    const lines = newCellLines[cellIndex];\n    const cell = table.getCell(r, c);\n    const body = cell.body;\n    body.paragraphs.load(\"items\");\n    await context.sync();\n\n    const paragraph = body.paragraphs.items[0];\n    const newText = lines.join(\"\\u000b\");\n    // Replace through the paragraph's own range so the existing run\n    // formatting (sz=32) is preserved instead of being dropped.\n    paragraph.getRange().insertText(newText, \"Replace\");\n    await context.sync();\n\n    cellIndex++;\n  }\n}\n", "ps1": "# Lattice multiplication exercises: update the numbers in each table cell.\n# Each cell holds a single run whose text is split into 5 segments joined by\n# manual line breaks (Chr(11) <-> <w:br/>):\n#   \"AB x CD\", \"  C    D\", \"  ----\", \"X|    |\", \"Y|    |\"\n# Table is 5 rows x 3 columns (row-major order below). Writing straight to\n# Cell.Range.Text keeps the cell's own run formatting (font size) intact.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$lineBreak = [char]11\n\n$newCells = @(\n    @(\"74 x 25\", \"  2    5\", \"  ----\", \"7|    |\", \"4|    |\"),\n    @(\"52 x 41\", \"  4    1\", \"  ----\", \"5|    |\", \"2|    |\"),\n    @(\"43 x 24\", \"  2    4\", \"  ----\", \"4|    |\", \"3|    |\"),\n    @(\"90 x 42\", \"  4    2\", \"  ----\", \"9|    |\", \"0|    |\"),\n    @(\"62 x 15\", \"  1    5\", \"  ----\", \"6|    |\", \"2|    |\"),\n    @(\"91 x 28\", \"  2    8\", \"  ----\", \"9|    |\", \"1|    |\"),\n    @(\"64 x 15\", \"  1    5\", \"  ----\", \"6|    |\", \"4|    |\"),\n    @(\"61 x 61\", \"  6    1\", \"  ----\", \"6|    |\", \"1|    |\"),\n    @(\"77 x 77\", \"  7    7\", \"  ----\", \"7|    |\", \"7|    |\"),\n    @(\"57 x 52\", \"  5    2\", \"  ----\", \"5|    |\", \"7|    |\"),\n    @(\"97 x 29\", \"  2    9\", \"  ----\", \"9|    |\", \"7|    |\"),\n    @(\"37 x 65\", \"  6    5\", \"  ----\", \"3|    |\", \"7|    |\"),\n    @(\"50 x 74\", \"  7    4\", \"  ----\", \"5|    |\", \"0|    |\"),\n    @(\"34 x 35\", \"  3    5\", \"  ----\", \"3|    |\", \"4|    |\"),\n    @(\"17 x 16\", \"  1    6\", \"  ----\", \"1|    |\", \"7|    |\")\n)\n\n$rows = 5\n$cols = 3\n$i = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $lines = $newCells[$i]\n        $newText = [string]::Join($lineBreak, $lines)\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $newText\n        $i++\n    }\n}\n"}
